$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.06742302838872502
$ws.Range("J2").Value = 0.06742302838872503
$ws.Range("M2").Value = 38.745275
$ws.Range("N2").Value = 116.235825
$ws.Range("O2").Value = 0.3160319337595895
$ws.Range("P2").Value = 0.3160319337595895
$ws.Range("Q2").Value = 5.623941241708333
$ws.Range("R2").Value = 50.61547117537499
$ws.Range("S2").Value = 0.02130783004161647
$ws.Range("T2").Value = 0.02130783004161647
$ws.Range("I3").Value = 0.06742302838872502
$ws.Range("J3").Value = 0.06742302838872503
$ws.Range("O3").Value = 0.4383510712400526
$ws.Range("P3").Value = 0.4383510712400526
$ws.Range("Q3").Value = 7.800669503763888
$ws.Range("R3").Value = 70.20602553387499
$ws.Range("S3").Value = 0.02955495672044609
$ws.Range("T3").Value = 0.0295549567204461
$ws.Range("I4").Value = 0.06742302838872502
$ws.Range("J4").Value = 0.06742302838872503
$ws.Range("M4").Value = 17.38482166666667
$ws.Range("N4").Value = 52.154465
$ws.Range("O4").Value = 0.1418020341675798
$ws.Range("P4").Value = 0.1418020341675798
$ws.Range("Q4").Value = 2.523435839619444
$ws.Range("R4").Value = 22.710922556575
$ws.Range("S4").Value = 0.00956072257525969
$ws.Range("T4").Value = 0.009560722575259693
$ws.Range("I5").Value = 0.06742302838872502
$ws.Range("J5").Value = 0.06742302838872503
$ws.Range("M5").Value = 12.72763533333333
$ws.Range("N5").Value = 38.182906
$ws.Range("O5").Value = 0.103814960832778
$ws.Range("P5").Value = 0.103814960832778
$ws.Range("Q5").Value = 1.847437481358889
$ws.Range("R5").Value = 16.62693733223
$ws.Range("S5").Value = 0.006999519051402765
$ws.Range("T5").Value = 0.006999519051402767
$ws.Range("G6").Value = 0.3560033333333333
$ws.Range("I6").Value = 0.1653637426357309
$ws.Range("J6").Value = 0.1653637426357309
$ws.Range("M6").Value = 38.745275
$ws.Range("N6").Value = 116.235825
$ws.Range("O6").Value = 0.3160319337595895
$ws.Range("P6").Value = 0.3160319337595895
$ws.Range("Q6").Value = 13.79344705091667
$ws.Range("R6").Value = 124.14102345825
$ws.Range("S6").Value = 0.05226022335889313
$ws.Range("T6").Value = 0.05226022335889313
$ws.Range("G7").Value = 0.3560033333333333
$ws.Range("I7").Value = 0.1653637426357309
$ws.Range("J7").Value = 0.1653637426357309
$ws.Range("O7").Value = 0.4383510712400526
$ws.Range("P7").Value = 0.4383510712400526
$ws.Range("S7").Value = 0.07248737372863701
$ws.Range("T7").Value = 0.07248737372863702
$ws.Range("G8").Value = 0.3560033333333333
$ws.Range("I8").Value = 0.1653637426357309
$ws.Range("J8").Value = 0.1653637426357309
$ws.Range("M8").Value = 17.38482166666667
$ws.Range("N8").Value = 52.154465
$ws.Range("O8").Value = 0.1418020341675798
$ws.Range("P8").Value = 0.1418020341675798
$ws.Range("Q8").Value = 6.189054462738889
$ws.Range("R8").Value = 55.70149016465
$ws.Range("S8").Value = 0.02344891508331079
$ws.Range("T8").Value = 0.0234489150833108
$ws.Range("G9").Value = 0.3560033333333333
$ws.Range("I9").Value = 0.1653637426357309
$ws.Range("J9").Value = 0.1653637426357309
$ws.Range("M9").Value = 12.72763533333333
$ws.Range("N9").Value = 38.182906
$ws.Range("O9").Value = 0.103814960832778
$ws.Range("P9").Value = 0.103814960832778
$ws.Range("Q9").Value = 4.531080604117777
$ws.Range("R9").Value = 40.77972543706
$ws.Range("S9").Value = 0.01716723046488998
$ws.Range("T9").Value = 0.01716723046488999
$ws.Range("G10").Value = 1.651695
$ws.Range("H10").Value = 4.955085
$ws.Range("I10").Value = 0.767213228975544
$ws.Range("J10").Value = 0.7672132289755441
$ws.Range("M10").Value = 38.745275
$ws.Range("N10").Value = 116.235825
$ws.Range("O10").Value = 0.3160319337595895
$ws.Range("P10").Value = 0.3160319337595895
$ws.Range("Q10").Value = 63.995376991125
$ws.Range("R10").Value = 575.9583929201249
$ws.Range("S10").Value = 0.2424638803590799
$ws.Range("T10").Value = 0.2424638803590799
$ws.Range("G11").Value = 1.651695
$ws.Range("H11").Value = 4.955085
$ws.Range("I11").Value = 0.767213228975544
$ws.Range("J11").Value = 0.7672132289755441
$ws.Range("O11").Value = 0.4383510712400526
$ws.Range("P11").Value = 0.4383510712400526
$ws.Range("Q11").Value = 88.764580606625
$ws.Range("R11").Value = 798.8812254596249
$ws.Range("S11").Value = 0.3363087407909695
$ws.Range("T11").Value = 0.3363087407909695
$ws.Range("G12").Value = 1.651695
$ws.Range("H12").Value = 4.955085
$ws.Range("I12").Value = 0.767213228975544
$ws.Range("J12").Value = 0.7672132289755441
$ws.Range("M12").Value = 17.38482166666667
$ws.Range("N12").Value = 52.154465
$ws.Range("O12").Value = 0.1418020341675798
$ws.Range("P12").Value = 0.1418020341675798
$ws.Range("Q12").Value = 28.714423022725
$ws.Range("R12").Value = 258.429807204525
$ws.Range("S12").Value = 0.1087923965090093
$ws.Range("T12").Value = 0.1087923965090094
$ws.Range("G13").Value = 1.651695
$ws.Range("H13").Value = 4.955085
$ws.Range("I13").Value = 0.767213228975544
$ws.Range("J13").Value = 0.7672132289755441
$ws.Range("M13").Value = 12.72763533333333
$ws.Range("N13").Value = 38.182906
$ws.Range("O13").Value = 0.103814960832778
$ws.Range("P13").Value = 0.103814960832778
$ws.Range("Q13").Value = 21.02217164189
$ws.Range("R13").Value = 189.19954477701
$ws.Range("S13").Value = 0.07964821131648522
$ws.Range("T13").Value = 0.07964821131648522
